$d = $word.ActiveDocument
$t = $d.Tables.Item(2)
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ---------------------------------------------------------------------------
# 1. New feature: Configuration -> Regional Country
#    Insert a new row "Company Information | Global | Regional country"
#    right after the existing "Company Information | Global | Language code"
#    row (and right before the following blank spacer row).
# ---------------------------------------------------------------------------
$langCodeRow = $null
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $row = $t.Rows.Item($r)
    $c1 = $row.Cells.Item(1).Range.Text
    $c3 = $row.Cells.Item(3).Range.Text
    if ($c1 -like "Company Information*" -and $c3 -like "Language*code*") {
        $langCodeRow = $r
        break
    }
}

$newRow = $t.Rows.Add($t.Rows.Item($langCodeRow + 1))
$newRow.Cells.Item(1).Range.Text = "Company Information"
$newRow.Cells.Item(2).Range.Text = "Global"
$newRow.Cells.Item(3).Range.Text = "Regional country"

# ---------------------------------------------------------------------------
# Helper: move the <w:lastRenderedPageBreak/> marker that sits in front of the
# first run of a table cell's paragraph from one row to another, preserving
# the paragraph/run formatting (fr-FR language mark).
# ---------------------------------------------------------------------------
function Set-FirstCellPageBreak($row, [bool]$withBreak) {
    $cell = $row.Cells.Item(1)
    $text = $cell.Range.Text
    $text = $text -replace "`r", ""
    $text = $text -replace "`a", ""

    $breakTag = ""
    if ($withBreak) { $breakTag = "<w:lastRenderedPageBreak/>" }

    $xml = "<w:p $wns><w:pPr><w:rPr><w:lang w:val='fr-FR'/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val='fr-FR'/></w:rPr>$breakTag<w:t xml:space='preserve'>$text</w:t></w:r></w:p>"

    $para = $cell.Range.Paragraphs.Item(1)
    $rng = $para.Range.Duplicate
    $rng.InsertXML($xml) | Out-Null
}

# ---------------------------------------------------------------------------
# 2. In the "TST parser tools" table block, the rendered page break marks
#    moved from the "cost_identifier" row onto the "Ticket identifier" row
#    (a natural side effect of the new row pushing content down a line).
# ---------------------------------------------------------------------------
$ticketIdRow = $null
$costIdRow = $null
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $row = $t.Rows.Item($r)
    $c1 = $row.Cells.Item(1).Range.Text
    $c3 = $row.Cells.Item(3).Range.Text
    if ($c1 -like "TST parser tools*") {
        if ($c3 -like "Ticket identifier*") { $ticketIdRow = $r }
        if ($c3 -like "cost_identifier*") { $costIdRow = $r }
    }
}

if ($ticketIdRow -ne $null) {
    Set-FirstCellPageBreak $t.Rows.Item($ticketIdRow) $true
}
if ($costIdRow -ne $null) {
    Set-FirstCellPageBreak $t.Rows.Item($costIdRow) $false
}

# ---------------------------------------------------------------------------
# 3. Likewise in the "Zenith Receipt parser tools" table block, the marker
#    moved from the "Ticket payment part" row onto the "Current travel
#    agency identifier" row above it.
# ---------------------------------------------------------------------------
$travelAgencyRow = $null
$ticketPaymentRow = $null
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $row = $t.Rows.Item($r)
    $c1 = $row.Cells.Item(1).Range.Text
    $c3 = $row.Cells.Item(3).Range.Text
    if ($c1 -like "Zenith Receipt parser tools*") {
        if ($c3 -like "Current travel agency identifier*") { $travelAgencyRow = $r }
        if ($c3 -like "Ticket payment part*") { $ticketPaymentRow = $r }
    }
}

if ($travelAgencyRow -ne $null) {
    Set-FirstCellPageBreak $t.Rows.Item($travelAgencyRow) $true
}
if ($ticketPaymentRow -ne $null) {
    Set-FirstCellPageBreak $t.Rows.Item($ticketPaymentRow) $false
}

Write-Output "done"
